# Running all the suites: flip every module's Runmode to "Y" so the whole
# suite executes, and leave the selection where the user finished editing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2:C6").Value = "Y"

$ws.Range("G7").Select()
